# Daily update at 8 AM UTC
# Appends the next day's row of data, and restores the "normal" date
# number format on the row that was previously the last row (A58),
# moving the special "last row" date format onto the new last row (A59).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous last row (58) had a distinct date-only number format
# applied because it was the final row. Now that a new row follows it,
# restore it to the regular date-time number format used by all other
# data rows.
$ws.Range("A58").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 59.
$ws.Range("A59").Value = 45799
$ws.Range("B59").Value = 245
$ws.Range("C59").Value = 256
$ws.Range("D59").Value = 246

# The newest row becomes the new "last row" and gets the distinct
# date-only number format.
$ws.Range("A59").NumberFormat = "YYYY-MM-DD"
